# Update the "quiz" marksheet for 1401ME56: correct marks / total marks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Marking row (B11): corrected-answer marks value
$ws.Range("B11").Value = 5

# Total row (B12): total score value
$ws.Range("B12").Value = 90

# Total row (E12): "correct/total" summary text
$ws.Range("E12").Value = "90/140"
